$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new account row (004526450 / MSD / 200000) right after the
#     header, i.e. as the new row 2 (pushes the existing data down by one). ---
$ws.Rows("2:2").Insert()

# Write the two text columns via a text-producing formula first (so the
# numeric-looking account number keeps its leading zeros / isn't coerced
# to a number), then flatten formula -> static value with a self
# copy/paste-special so no stray number-format style gets introduced.
$ws.Cells.Item(2, 1).Formula = "=""004526450"""
$ws.Cells.Item(2, 2).Formula = "=""MSD"""
$ws.Cells.Item(2, 3).Value = 200000

$ws.Range("A2:C2").Copy()
$ws.Range("A2:C2").PasteSpecial(-4163)

# --- Remove the duplicate "004752494 / SERGIO / 6623.66" row. It was
#     originally row 7; after inserting the row above it is now row 8. ---
$ws.Rows("8:8").Delete()
